# Snapshot and nicolas comments
#
# 1. Remove column M ("Category"/"ECE") from every sheet - it was a
#    constant column that added no information.
# 2. Re-word the severity headers in row 1 (columns E:L) to be more
#    descriptive, e.g. "% 1-2" -> "% severity levels 1-2".

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Drop the whole "Category" column (M) - shifts nothing else around
    # since it was the last column.
    $ws.Columns.Item(13).Delete()

    # Clarify the header labels for the severity-level columns.
    $ws.Range("E1").Value = "% severity levels 1-2"
    $ws.Range("F1").Value = "# severity levels 1-2"
    $ws.Range("G1").Value = "% severity level 3"
    $ws.Range("H1").Value = "# severity level 3"
    $ws.Range("I1").Value = "% severity level 4"
    $ws.Range("J1").Value = "# severity level 4"
    $ws.Range("K1").Value = "% severity level 5"
    $ws.Range("L1").Value = "# severity level 5"
}
